# Updates the cryptos list (commit: "Updated cryptos list ... with GitHub Actions").
# Source diff touches columns D (Price) and E (Volume(1h)) for most rows, and
# additionally swaps the Gas/HuobiToken rows B/C/D/E content (rows 48-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value. Values are written as literal text (matching the workbooks
# original inline-string storage) even when they look numeric, e.g. "245.71" or
# "8.50", so trailing zeros / multi-dot "prices" like "35.311.10" are preserved.
$updates = [ordered]@{
    "D2" = "35.311.10"
    "E2" = "  +0.41%  "
    "D3" = "1.883.88"
    "E3" = "  -0.89%  "
    "E4" = "  -0.68%  "
    "D5" = "245.71"
    "D6" = "0.690"
    "E6" = "  -0.95%  "
    "E7" = "  -0.72%  "
    "D8" = "43.37"
    "E8" = "  +5.96%  "
    "D9" = "0.355"
    "E9" = "  -1.30%  "
    "D10" = "53.55"
    "E10" = "  +1.28%  "
    "E11" = "  -1.73%  "
    "E12" = "  -1.21%  "
    "D13" = "13.24"
    "E13" = "  +1.85%  "
    "D14" = "2.158.96"
    "E14" = "  -0.79%  "
    "E15" = "  +2.55%  "
    "D16" = "4.89"
    "E16" = "  -1.61%  "
    "D17" = "1.901.36"
    "E17" = "  -0.20%  "
    "D18" = "35.415.29"
    "E18" = "  +0.73%  "
    "E19" = "  -1.01%  "
    "E20" = "  -1.59%  "
    "D21" = "244.47"
    "E22" = "  -1.65%  "
    "E23" = "  -2.13%  "
    "D24" = "2.64"
    "E24" = "  +9.03%  "
    "E25" = "  -0.74%  "
    "D26" = "2.13"
    "E26" = "  -7.26%  "
    "D27" = "165.76"
    "E27" = "  -0.67%  "
    "D28" = "8.50"
    "D29" = "18.27"
    "E29" = "  -1.17%  "
    "E30" = "  -2.13%  "
    "D31" = "4.128.46"
    "E31" = "  +0.01%  "
    "D32" = "1.72"
    "E32" = "  +9.42%  "
    "E33" = "  -1.25%  "
    "E34" = "  -4.45%  "
    "E35" = "  -1.30%  "
    "E36" = "  -0.76%  "
    "E37" = "  -11.45%  "
    "D38" = "0.849"
    "E38" = "  -0.11%  "
    "D39" = "1.95"
    "E39" = "  -2.76%  "
    "E40" = "  +7.13%  "
    "E41" = "  +2.95%  "
    "D42" = "17.24"
    "E42" = "  -0.10%  "
    "D43" = "96.56"
    "E43" = "  -5.31%  "
    "E44" = "  -2.30%  "
    "D45" = "1.299.35"
    "E45" = "  -1.63%  "
    "E46" = "  -4.91%  "
    "D47" = "0.0798"
    "E47" = "  +7.90%  "
    "B48" = "HuobiToken"
    "C48" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D48" = "2.37"
    "E48" = "  -2.31%  "
    "B49" = "Gas"
    "C49" = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
    "D49" = "12.27"
    "E49" = "  +3.51%  "
    "E50" = "  -0.78%  "
    "E51" = "  -5.38%  "
}

# Cells whose new value would otherwise be auto-coerced to a Number by the
# COM Value setter (plain decimals like "245.71"); force text formatting for
# those specific cells, then restore the default "Normal" style so the rest of
# the cell formatting (borders/alignment/etc.) is unaffected.
$forceText = @(
    "D5", "D6", "D8", "D9", "D10", "D13", "D16", "D21", "D24", "D26", "D27", "D28", "D29", "D32", "D38", "D39", "D42", "D43", "D47", "D48", "D49"
)

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    if ($forceText -contains $cellRef) {
        $range.NumberFormat = "@"
        $range.Value = $updates[$cellRef]
        $range.Style = "Normal"
    } else {
        $range.Value = $updates[$cellRef]
    }
}
